# Applies the "final grade view" edit: updates the rollno (column G) and
# name (column H) values for rows 2-16 on Sheet1, shifting each row's
# original values up by one row (row 2's original values wrap around to
# row 16), normalizing the "lci2021051" roll number to uppercase.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rollno = @{
    2  = "LCI2021009"
    3  = "LCI2021001"
    4  = "LCI2021017"
    5  = "LCI2021024"
    6  = "LCI2021002"
    7  = "LCI2021007"
    8  = "LCI2021018"
    9  = "LCI2021010"
    10 = "LCI2021057"
    11 = "LCI2021011"
    12 = "LCI2021027"
    13 = "LCI2021035"
    14 = "LCI2021043"
    15 = "LCI2021045"
    16 = "LCI2021051"
}

$name = @{
    2  = "Gaurav Kabra"
    3  = "Atharv Tiwari"
    4  = "Prateek Parmar"
    5  = "Rahul Jamwal"
    6  = "Harsh Golchha"
    7  = "Sameep Aher"
    8  = "Vidisha Agarwal"
    9  = "Samarth Sharma"
    10 = "Rhysha Kachari"
    11 = "Chakradhar Reddy"
    12 = "Yatharth Jain"
    13 = "Saarthak Verma"
    14 = "Bhavya Choudhary"
    15 = "Advit Mahale"
    16 = "Shruti Gajbhiye"
}

foreach ($row in 2..16) {
    $ws.Cells.Item($row, 7).Value = $rollno[$row]
    $ws.Cells.Item($row, 8).Value = $name[$row]
}
